$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "B2m"
$ws.Cells.Item(2,3).Value = "Gm11127"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1790.094238333333
$ws.Cells.Item(2,8).Value = 5370.282715
$ws.Cells.Item(2,9).Value = 0.1978827026976269
$ws.Cells.Item(2,10).Value = 0.1978827026976269
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.359983
$ws.Cells.Item(2,14).Value = 1.079949
$ws.Cells.Item(2,15).Value = 0.9497806167005701
$ws.Cells.Item(2,16).Value = 0.9497806167005701
$ws.Cells.Item(2,17).Value = 644.4034941979484
$ws.Cells.Item(2,18).Value = 5799.631447781536
$ws.Cells.Item(2,19).Value = 0.1879451554025277
$ws.Cells.Item(2,20).Value = 0.1879451554025277

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "B2m"
$ws.Cells.Item(3,3).Value = "Gm11127"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1790.094238333333
$ws.Cells.Item(3,8).Value = 5370.282715
$ws.Cells.Item(3,9).Value = 0.1978827026976269
$ws.Cells.Item(3,10).Value = 0.1978827026976269
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.019034
$ws.Cells.Item(3,14).Value = 0.057102
$ws.Cells.Item(3,15).Value = 0.05021938329942984
$ws.Cells.Item(3,16).Value = 0.05021938329942984
$ws.Cells.Item(3,17).Value = 34.07265373243666
$ws.Cells.Item(3,18).Value = 306.65388359193
$ws.Cells.Item(3,19).Value = 0.009937547295099245
$ws.Cells.Item(3,20).Value = 0.009937547295099247

$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "B2m"
$ws.Cells.Item(4,3).Value = "Gm11127"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1072.362365666667
$ws.Cells.Item(4,8).Value = 3217.087097
$ws.Cells.Item(4,9).Value = 0.1185423418752029
$ws.Cells.Item(4,10).Value = 0.1185423418752029
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.359983
$ws.Cells.Item(4,14).Value = 1.079949
$ws.Cells.Item(4,15).Value = 0.9497806167005701
$ws.Cells.Item(4,16).Value = 0.9497806167005701
$ws.Cells.Item(4,17).Value = 386.0322214797836
$ws.Cells.Item(4,18).Value = 3474.289993318053
$ws.Cells.Item(4,19).Value = 0.11258921857136
$ws.Cells.Item(4,20).Value = 0.11258921857136

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "B2m"
$ws.Cells.Item(5,3).Value = "Gm11127"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1072.362365666667
$ws.Cells.Item(5,8).Value = 3217.087097
$ws.Cells.Item(5,9).Value = 0.1185423418752029
$ws.Cells.Item(5,10).Value = 0.1185423418752029
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.019034
$ws.Cells.Item(5,14).Value = 0.057102
$ws.Cells.Item(5,15).Value = 0.05021938329942984
$ws.Cells.Item(5,16).Value = 0.05021938329942984
$ws.Cells.Item(5,17).Value = 20.41134526809933
$ws.Cells.Item(5,18).Value = 183.702107412894
$ws.Cells.Item(5,19).Value = 0.005953123303842865
$ws.Cells.Item(5,20).Value = 0.005953123303842866

$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "B2m"
$ws.Cells.Item(6,3).Value = "Gm11127"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 2676.771728333333
$ws.Cells.Item(6,8).Value = 8030.315184999999
$ws.Cells.Item(6,9).Value = 0.2958988486552321
$ws.Cells.Item(6,10).Value = 0.2958988486552321
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 0.6666666666666666
$ws.Cells.Item(6,13).Value = 0.359983
$ws.Cells.Item(6,14).Value = 1.079949
$ws.Cells.Item(6,15).Value = 0.9497806167005701
$ws.Cells.Item(6,16).Value = 0.9497806167005701
$ws.Cells.Item(6,17).Value = 963.5923170806183
$ws.Cells.Item(6,18).Value = 8672.330853725565
$ws.Cells.Item(6,19).Value = 0.281038990956755
$ws.Cells.Item(6,20).Value = 0.281038990956755

$ws.Cells.Item(7,1).Value = "M1"
$ws.Cells.Item(7,2).Value = "B2m"
$ws.Cells.Item(7,3).Value = "Gm11127"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2676.771728333333
$ws.Cells.Item(7,8).Value = 8030.315184999999
$ws.Cells.Item(7,9).Value = 0.2958988486552321
$ws.Cells.Item(7,10).Value = 0.2958988486552321
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 0.019034
$ws.Cells.Item(7,14).Value = 0.057102
$ws.Cells.Item(7,15).Value = 0.05021938329942984
$ws.Cells.Item(7,16).Value = 0.05021938329942984
$ws.Cells.Item(7,17).Value = 50.94967307709666
$ws.Cells.Item(7,18).Value = 458.54705769387
$ws.Cells.Item(7,19).Value = 0.01485985769847708
$ws.Cells.Item(7,20).Value = 0.01485985769847708

$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "B2m"
$ws.Cells.Item(8,3).Value = "Gm11127"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 3409.239257666667
$ws.Cells.Item(8,8).Value = 10227.717773
$ws.Cells.Item(8,9).Value = 0.3768681357681173
$ws.Cells.Item(8,10).Value = 0.3768681357681174
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.359983
$ws.Cells.Item(8,14).Value = 1.079949
$ws.Cells.Item(8,15).Value = 0.9497806167005701
$ws.Cells.Item(8,16).Value = 0.9497806167005701
$ws.Cells.Item(8,17).Value = 1227.26817569262
$ws.Cells.Item(8,18).Value = 11045.41358123358
$ws.Cells.Item(8,19).Value = 0.3579420504046367
$ws.Cells.Item(8,20).Value = 0.3579420504046367

$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "B2m"
$ws.Cells.Item(9,3).Value = "Gm11127"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 3409.239257666667
$ws.Cells.Item(9,8).Value = 10227.717773
$ws.Cells.Item(9,9).Value = 0.3768681357681173
$ws.Cells.Item(9,10).Value = 0.3768681357681174
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 0.019034
$ws.Cells.Item(9,14).Value = 0.057102
$ws.Cells.Item(9,15).Value = 0.05021938329942984
$ws.Cells.Item(9,16).Value = 0.05021938329942984
$ws.Cells.Item(9,17).Value = 64.89146003042733
$ws.Cells.Item(9,18).Value = 584.023140273846
$ws.Cells.Item(9,19).Value = 0.01892608536348065
$ws.Cells.Item(9,20).Value = 0.01892608536348065

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "B2m"
$ws.Cells.Item(10,3).Value = "Gm11127"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 97.77148966666668
$ws.Cells.Item(10,8).Value = 293.314469
$ws.Cells.Item(10,9).Value = 0.01080797100382067
$ws.Cells.Item(10,10).Value = 0.01080797100382067
$ws.Cells.Item(10,11).Value = 2
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.359983
$ws.Cells.Item(10,14).Value = 1.079949
$ws.Cells.Item(10,15).Value = 0.9497806167005701
$ws.Cells.Item(10,16).Value = 0.9497806167005701
$ws.Cells.Item(10,17).Value = 35.19607416467567
$ws.Cells.Item(10,18).Value = 316.7646674820811
$ws.Cells.Item(10,19).Value = 0.01026520136529067
$ws.Cells.Item(10,20).Value = 0.01026520136529067

$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "B2m"
$ws.Cells.Item(11,3).Value = "Gm11127"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 97.77148966666668
$ws.Cells.Item(11,8).Value = 293.314469
$ws.Cells.Item(11,9).Value = 0.01080797100382067
$ws.Cells.Item(11,10).Value = 0.01080797100382067
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.019034
$ws.Cells.Item(11,14).Value = 0.057102
$ws.Cells.Item(11,15).Value = 0.05021938329942984
$ws.Cells.Item(11,16).Value = 0.05021938329942984
$ws.Cells.Item(11,17).Value = 1.860982534315333
$ws.Cells.Item(11,18).Value = 16.748842808838
$ws.Cells.Item(11,19).Value = 0.0005427696385299935
$ws.Cells.Item(11,20).Value = 0.0005427696385299935
